$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $wsHoja1.Range("A1")
$text = $cell.Value()
$text = $text.Replace("1000 Bs = 9.94 = 42533.04 pesos", "1000 Bs = 9.91 = 42200.2 pesos")
$text = $text.Replace("42533.04 pesos = 9.91 = 954.74 Bs", "42200.2 pesos = 9.83 = 953.51 Bs")
$cell.Value = $text

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 100.9
$wsTasas.Range("O10").Value = 4258
$wsTasas.Range("N12").Value = 4292.98
$wsTasas.Range("O12").Value = 97
